$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 65; this shifts the existing rows 65-171 down to 66-172.
$ws.Rows.Item(65).Insert()

# Populate the newly inserted row 65 with its data (same structure as the other
# data rows, with the new record's values).
$ws.Cells.Item(65, 1).Value = 5
$ws.Cells.Item(65, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(65, 3).Value = "Maule"
$ws.Cells.Item(65, 4).Value = 45272
$ws.Cells.Item(65, 5).Value = 7
$ws.Cells.Item(65, 6).Value = 100112022
$ws.Cells.Item(65, 7).Value = "Arveja Verde"
$ws.Cells.Item(65, 8).Value = "Sin especificar"
$ws.Cells.Item(65, 9).Value = "Primera"
$ws.Cells.Item(65, 10).Value = 300
$ws.Cells.Item(65, 11).Value = 23000
$ws.Cells.Item(65, 12).Value = 23000
$ws.Cells.Item(65, 13).Value = 23000
$ws.Cells.Item(65, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(65, 15).Value = "Región del Maule"
$ws.Cells.Item(65, 16).Value = 920
$ws.Cells.Item(65, 17).Value = 25
$ws.Cells.Item(65, 18).Value = "Hortaliza"
